# Added buttons to flavors page
$wb = $excel.ActiveWorkbook

# Rename the "Gluten-Free Options" sheet to "Gluten-Free"
$glutenFreeSheet = $wb.Worksheets.Item("Gluten-Free Options")
$glutenFreeSheet.Name = "Gluten-Free"

# Make the Gluten-Free sheet the active tab (was previously Classic)
$glutenFreeSheet.Activate()
